$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.726787269115448
$ws.Range("B1").Value = 2.112131118774414
$ws.Range("C1").Value = 2.545148611068726
$ws.Range("D1").Value = 3.025104999542236
$ws.Range("E1").Value = 1.169565081596375
